$d = $word.ActiveDocument

$find = "Anomaly Detection (Deep One-Class Learning), Evolutionary Algorithm, Graph Theory and "
$replace = "Anomaly Detection (Deep One-Class Learning), Optimal Control and Estimation, Evolutionary Algorithm, Graph Theory and "

$d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
